$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B8").Value = "2026-01-16T13:49:34+00:00"
$ws1.Range("B12").Value = "Acte"

$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("M2").Value = "Acte"
$ws2.Range("L13").Value = "Auteur"
$ws2.Range("M13").Value = "Auteur"
$ws2.Range("L14").Value = "Informateur"
$ws2.Range("M14").Value = "Informateur"
$ws2.Range("L15").Value = "Participant"
$ws2.Range("M15").Value = "Participant"
